{"js": "// The screw size used to mount the motor was changed from \"4-32\" to\n// \"4-40\" in two places in the document: the parts list bullet\n// (\"4-32 x \u00bc screw\") and the build-instructions paragraph describing\n// how the motor is mounted (\"...screw (4-32 x \u00bc)...\").\nconst body = context.document.body;\nconst results = body.search(\"4-32\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < results.items.length; i++) {\n  results.items[i].insertText(\"4-40\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "# The screw size used to mount the motor was changed from \"4-32\" to\n# \"4-40\" in two places in the document: the parts list bullet\n# (\"4-32 x \u00bc screw\") and the build-instructions paragraph describing\n# how the motor is mounted (\"...screw (4-32 x \u00bc)...\").\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \"4-32\"\n$find.Replacement.Text = \"4-40\"\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$find.Execute($find.Text, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $find.Replacement.Text, $wdReplaceAll)\n"}
